$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I20").Value = -0.05808045239458798
$ws.Range("J20").Value = 0.2399996375480334
$ws.Range("K20").Value = 0.237099825084014
$ws.Range("L20").Value = 2.554812062321973
